$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (same set applied to every data row, columns B:Q)
$newValues = @(
    [double]"0.9999674344715328",
    [double]"0.9989400190852257",
    [double]"0.999954822450698",
    [double]"0.9999963708365882",
    [double]"0.9999770133602924",
    [double]"3.039849136261971e-05",
    [double]"0.0009894456561571146",
    [double]"4.373235006023419e-05",
    [double]"4.008672657845657e-06",
    [double]"2.387051135903992e-05",
    [double]"0.000349036811621329",
    [double]"0.005513482689065026",
    [double]"1.000060120975632",
    [double]"0.005748203061621094",
    [double]"94.80223515402628",
    [double]"139.9006406741497"
)

for ($row = 2; $row -le 26; $row++) {
    for ($i = 0; $i -lt $newValues.Length; $i++) {
        $col = 2 + $i   # Column B = 2 ... Column Q = 17
        $ws.Cells.Item($row, $col).Value = $newValues[$i]
    }
}
